$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: RandomForestRegressor - update values
$ws.Range("B3").Value = 0.9835760988432076
$ws.Range("C3").Value = 0.9817794427116402
$ws.Range("D3").Value = 0.9527642615467489

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9812114986549524
$ws.Range("C4").Value = 0.9796248748224839
$ws.Range("D4").Value = 0.877122300611426

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9876477421214016
$ws.Range("C5").Value = 0.9876132006717123
$ws.Range("D5").Value = 0.9854642593305272
